# Re-sort the "estado de cuenta" detail table (rows 16-50 on Hoja1) so that it is
# grouped by Periodo Mora (ascending: 2408 -> 2502) with the five workers repeated
# inside each period, instead of grouped by worker (descending period inside each
# worker) as it was before. The (worker, period) -> (Valor Mora) values themselves
# are unchanged; only the row order / grouping changes. Also adds the "parte 1" of
# the new estado de cuenta data per the commit message (same underlying dataset,
# reorganized).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Workers in the same order they already appear in the sheet (doc number, name).
$workers = @(
    @{ Doc = "1047390509"; Name = "JAVIER JESUS VARGAS PEREZ" },
    @{ Doc = "1127584034"; Name = "ALONSO JOSE VELASCO CARRILLO" },
    @{ Doc = "1103117470"; Name = "JUAN CAMILO LOPEZ RIOS" },
    @{ Doc = "1085038750"; Name = "RAUL RUIZ RAMOS" },
    @{ Doc = "1052071317"; Name = "OLGA ISABEL VASQUEZ TEHERAN" }
)

# Periods, now ascending (oldest first).
$periods = @("2408", "2409", "2410", "2411", "2412", "2501", "2502")

# Valor Mora depends only on the period.
function Get-ValorMora($periodo) {
    if ($periodo -eq "2502") { return 32933 }
    if ($periodo -eq "2408") { return 45066 }
    return 52000
}

$row = 16
foreach ($periodo in $periods) {
    foreach ($worker in $workers) {
        $ws.Cells.Item($row, 3).Value = $worker.Doc
        $ws.Cells.Item($row, 4).Value = $worker.Name
        $ws.Cells.Item($row, 5).Value = $periodo
        $ws.Cells.Item($row, 6).Value = Get-ValorMora $periodo
        $row++
    }
}
